# Fruta / hortaliza, semanal
# Insert a new weekly observation row (row 191) for the Mango subset, shifting
# the former last row down to 192, and refresh the values of the rows that
# the weekly reload touched (188-192).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new weekly record: everything from the old row 191 on
# shifts down by one.
$ws.Rows(191).Insert()

function Set-MangoRow($Row, $Date, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Calidad, $Unidad, $PrecioKg, $KgUnidad) {
    $ws.Cells.Item($Row, 1).Value = 1
    $ws.Cells.Item($Row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($Row, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($Row, 4).Value = $Date
    $ws.Cells.Item($Row, 5).Value = 15
    $ws.Cells.Item($Row, 6).Value = "Fruta"
    $ws.Cells.Item($Row, 7).Value = 100108
    $ws.Cells.Item($Row, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($Row, 9).Value = 100108002
    $ws.Cells.Item($Row, 10).Value = "Mango"
    $ws.Cells.Item($Row, 11).Value = "Sin especificar"
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = "Perú"
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

# Row 188: refreshed with this week's new report (date + prices bumped, unit
# switched from "caja" to "bandeja").
Set-MangoRow 188 44939 456 5000 5500 5250 "Especial" "$/bandeja 4 kilos" 1312 4

# Row 189: only the "Calidad" label was corrected.
$ws.Cells.Item(189, 12).Value = "Especial"

# Row 190: now carries what used to be the "caja" / Primera record.
Set-MangoRow 190 44249 450 4500 5000 4750 "Primera" "$/caja 4 kilos" 4750 1

# Row 191 (newly inserted): carries what used to be row 190's "bandeja" /
# Especial record.
Set-MangoRow 191 44356 400 3500 4000 3750 "Especial" "$/bandeja 4 kilos" 938 4

# Row 192: the former row 191, now shifted down one position, values unchanged.
Set-MangoRow 192 44335 456 4500 5000 4750 "Especial" "$/bandeja 4 kilos" 1188 4
